# This edit re-shuffles the per-row "observation" data (date, volume, prices,
# unit, origin, $/Kg, Kg/unit) among the existing data rows (2-26) of the
# sheet. The dimension columns (Mercado, Region, Producto, etc.) are
# identical across all rows, so only columns D and M:T need to move.
#
# Mapping is: destinationRow -> sourceRow (the row whose D/M:T values should
# end up at destinationRow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 21
    3  = 19
    4  = 7
    5  = 8
    6  = 20
    7  = 26
    8  = 17
    9  = 3
    10 = 22
    11 = 9
    12 = 16
    13 = 23
    14 = 10
    15 = 6
    16 = 15
    17 = 5
    18 = 11
    19 = 4
    20 = 12
    21 = 2
    22 = 18
    23 = 14
    24 = 25
    25 = 13
    26 = 24
}

$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

# First snapshot all the source rows' current values before we start writing,
# so that writes to earlier rows don't clobber data that is still needed as
# a source for a later row.
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
